# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# The sheet lists 9 workers in arrears ("cartera") for period 2508.
# This edit adds a second block of the same 9 workers for period 2509
# right below the existing block (duplicating rows 16:24 into new rows
# 25:33), bumps the summary totals (VALOR MORA / Cant. Periodos)
# accordingly, and keeps the table's banded-row formatting consistent
# (only the very last data row keeps the heavier "closing" bottom
# border; every other data row - including the old last row, which is
# no longer last - uses the regular interior-row border). The
# "Periodo Mora" column is also centered to match the rest of the
# numeric-looking columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the 9 worker rows (16:24) into a new block (25:33) ---
# Insert 9 blank rows right after the existing block, shifting the
# footer ("___", "NOMBRE DEL REPRESENTANTE LEGAL", ...) down from
# rows 29:30 to rows 38:39.
$ws.Range("A25:J33").Insert()

# Copy the whole original block (values + formatting) into the new rows.
$ws.Range("B16:J24").Copy($ws.Range("B25:J33"))

# The new block is for period 2509 instead of 2508.
for ($r = 25; $r -le 33; $r++) {
    $ws.Range("E" + $r).Value2 = "2509"
}

# --- 2. Fix up row banding now that row 24 is no longer the last row ---
# Row 24 (old last row) must go back to the regular interior-row
# formatting (it copies cleanly from row 23, its new neighbour);
# row 33 (the new last row) already received the "closing" formatting
# from the Copy above, since it came from old row 24.
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Center the "Periodo Mora" column across the whole table body ---
$ws.Range("E16:E33").HorizontalAlignment = -4108

# --- 4. Update the summary figures for the added period ---
$ws.Range("E11").Value2 = 1010512
$ws.Range("F13").Value2 = 2
